$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp header
$ws.Range("A1").Value = "Datos actualizados a 30 de Junio de 2020 a las 09:41"

# Update country rows whose case counts (and in a few cases, ranking/order) changed
$ws.Cells.Item(7, 1).Value = "India"
$ws.Cells.Item(7, 2).Value = 568315
$ws.Cells.Item(7, 3).Value = 779
$ws.Cells.Item(7, 4).Value = 335577
$ws.Cells.Item(7, 5).Value = 215821
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 13
$ws.Cells.Item(7, 8).Value = 16917

$ws.Cells.Item(37, 1).Value = "Ucrania"
$ws.Cells.Item(37, 2).Value = 44334
$ws.Cells.Item(37, 3).Value = 706
$ws.Cells.Item(37, 4).Value = 19115
$ws.Cells.Item(37, 5).Value = 24060
$ws.Cells.Item(37, 6).Value = 0
$ws.Cells.Item(37, 7).Value = 12
$ws.Cells.Item(37, 8).Value = 1159

$ws.Cells.Item(38, 1).Value = "Singapur"
$ws.Cells.Item(38, 2).Value = 43907
$ws.Cells.Item(38, 3).Value = 246
$ws.Cells.Item(38, 4).Value = 37985
$ws.Cells.Item(38, 5).Value = 5896
$ws.Cells.Item(38, 6).Value = 0
$ws.Cells.Item(38, 7).Value = 0
$ws.Cells.Item(38, 8).Value = 26

$ws.Cells.Item(50, 1).Value = "Armenia"
$ws.Cells.Item(50, 2).Value = 25542
$ws.Cells.Item(50, 3).Value = 415
$ws.Cells.Item(50, 4).Value = 14048
$ws.Cells.Item(50, 5).Value = 11051
$ws.Cells.Item(50, 6).Value = 0
$ws.Cells.Item(50, 7).Value = 10
$ws.Cells.Item(50, 8).Value = 443

$ws.Cells.Item(51, 1).Value = "Irlanda"
$ws.Cells.Item(51, 2).Value = 25462
$ws.Cells.Item(51, 3).Value = 0
$ws.Cells.Item(51, 4).Value = 23364
$ws.Cells.Item(51, 5).Value = 363
$ws.Cells.Item(51, 6).Value = 0
$ws.Cells.Item(51, 7).Value = 0
$ws.Cells.Item(51, 8).Value = 1735

$ws.Cells.Item(52, 1).Value = "Nigeria"
$ws.Cells.Item(52, 2).Value = 25133
$ws.Cells.Item(52, 3).Value = 0
$ws.Cells.Item(52, 4).Value = 9402
$ws.Cells.Item(52, 5).Value = 15158
$ws.Cells.Item(52, 6).Value = 0
$ws.Cells.Item(52, 7).Value = 0
$ws.Cells.Item(52, 8).Value = 573

$ws.Cells.Item(69, 1).Value = "Chequia"
$ws.Cells.Item(69, 2).Value = 11809
$ws.Cells.Item(69, 3).Value = 4
$ws.Cells.Item(69, 4).Value = 7751
$ws.Cells.Item(69, 5).Value = 3709
$ws.Cells.Item(69, 6).Value = 0
$ws.Cells.Item(69, 7).Value = 1
$ws.Cells.Item(69, 8).Value = 349

$ws.Cells.Item(94, 1).Value = "Hungria"
$ws.Cells.Item(94, 2).Value = 4155
$ws.Cells.Item(94, 3).Value = 10
$ws.Cells.Item(94, 4).Value = 2692
$ws.Cells.Item(94, 5).Value = 878
$ws.Cells.Item(94, 6).Value = 0
$ws.Cells.Item(94, 7).Value = 0
$ws.Cells.Item(94, 8).Value = 585

$ws.Cells.Item(104, 1).Value = "Estado de Palestina"
$ws.Cells.Item(104, 2).Value = 2390
$ws.Cells.Item(104, 3).Value = 205
$ws.Cells.Item(104, 4).Value = 451
$ws.Cells.Item(104, 5).Value = 1934
$ws.Cells.Item(104, 6).Value = 0
$ws.Cells.Item(104, 7).Value = 0
$ws.Cells.Item(104, 8).Value = 5

$ws.Cells.Item(105, 1).Value = "Cuba"
$ws.Cells.Item(105, 2).Value = 2340
$ws.Cells.Item(105, 3).Value = 0
$ws.Cells.Item(105, 4).Value = 2211
$ws.Cells.Item(105, 5).Value = 43
$ws.Cells.Item(105, 6).Value = 0
$ws.Cells.Item(105, 7).Value = 0
$ws.Cells.Item(105, 8).Value = 86

$ws.Cells.Item(106, 1).Value = "Maldivas"
$ws.Cells.Item(106, 2).Value = 2337
$ws.Cells.Item(106, 3).Value = 0
$ws.Cells.Item(106, 4).Value = 1927
$ws.Cells.Item(106, 5).Value = 402
$ws.Cells.Item(106, 6).Value = 0
$ws.Cells.Item(106, 7).Value = 0
$ws.Cells.Item(106, 8).Value = 8

$ws.Cells.Item(107, 1).Value = "Paraguay"
$ws.Cells.Item(107, 2).Value = 2191
$ws.Cells.Item(107, 3).Value = 0
$ws.Cells.Item(107, 4).Value = 1080
$ws.Cells.Item(107, 5).Value = 1095
$ws.Cells.Item(107, 6).Value = 0
$ws.Cells.Item(107, 7).Value = 0
$ws.Cells.Item(107, 8).Value = 16

$ws.Cells.Item(111, 1).Value = "Sri Lanka"
$ws.Cells.Item(111, 2).Value = 2042
$ws.Cells.Item(111, 3).Value = 3
$ws.Cells.Item(111, 4).Value = 1711
$ws.Cells.Item(111, 5).Value = 320
$ws.Cells.Item(111, 6).Value = 0
$ws.Cells.Item(111, 7).Value = 0
$ws.Cells.Item(111, 8).Value = 11

$ws.Cells.Item(114, 1).Value = "Estonia"
$ws.Cells.Item(114, 2).Value = 1989
$ws.Cells.Item(114, 3).Value = 2
$ws.Cells.Item(114, 4).Value = 1829
$ws.Cells.Item(114, 5).Value = 91
$ws.Cells.Item(114, 6).Value = 0
$ws.Cells.Item(114, 7).Value = 0
$ws.Cells.Item(114, 8).Value = 69

$ws.Cells.Item(116, 1).Value = "Lituania"
$ws.Cells.Item(116, 2).Value = 1817
$ws.Cells.Item(116, 3).Value = 1
$ws.Cells.Item(116, 4).Value = 1515
$ws.Cells.Item(116, 5).Value = 224
$ws.Cells.Item(116, 6).Value = 0
$ws.Cells.Item(116, 7).Value = 0
$ws.Cells.Item(116, 8).Value = 78

$ws.Cells.Item(131, 1).Value = "Letonia"
$ws.Cells.Item(131, 2).Value = 1118
$ws.Cells.Item(131, 3).Value = 1
$ws.Cells.Item(131, 4).Value = 974
$ws.Cells.Item(131, 5).Value = 114
$ws.Cells.Item(131, 6).Value = 0
$ws.Cells.Item(131, 7).Value = 0
$ws.Cells.Item(131, 8).Value = 30

$ws.Cells.Item(138, 1).Value = "Georgia"
$ws.Cells.Item(138, 2).Value = 928
$ws.Cells.Item(138, 3).Value = 2
$ws.Cells.Item(138, 4).Value = 794
$ws.Cells.Item(138, 5).Value = 119
$ws.Cells.Item(138, 6).Value = 0
$ws.Cells.Item(138, 7).Value = 0
$ws.Cells.Item(138, 8).Value = 15

$ws.Cells.Item(157, 1).Value = "Taiwan"
$ws.Cells.Item(157, 2).Value = 447
$ws.Cells.Item(157, 3).Value = 0
$ws.Cells.Item(157, 4).Value = 437
$ws.Cells.Item(157, 5).Value = 3
$ws.Cells.Item(157, 6).Value = 0
$ws.Cells.Item(157, 7).Value = 0
$ws.Cells.Item(157, 8).Value = 7

$ws.Cells.Item(193, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(193, 2).Value = 41
$ws.Cells.Item(193, 3).Value = 0
$ws.Cells.Item(193, 4).Value = 11
$ws.Cells.Item(193, 5).Value = 28
$ws.Cells.Item(193, 6).Value = 0
$ws.Cells.Item(193, 7).Value = 1
$ws.Cells.Item(193, 8).Value = 2

$ws.Cells.Item(205, 1).Value = "Fiyi"
$ws.Cells.Item(205, 2).Value = 18
$ws.Cells.Item(205, 3).Value = 0
$ws.Cells.Item(205, 4).Value = 18
$ws.Cells.Item(205, 5).Value = 0
$ws.Cells.Item(205, 6).Value = 0
$ws.Cells.Item(205, 7).Value = 0
$ws.Cells.Item(205, 8).Value = 0

$ws.Cells.Item(206, 1).Value = "Dominica"
$ws.Cells.Item(206, 2).Value = 18
$ws.Cells.Item(206, 3).Value = 0
$ws.Cells.Item(206, 4).Value = 18
$ws.Cells.Item(206, 5).Value = 0
$ws.Cells.Item(206, 6).Value = 0
$ws.Cells.Item(206, 7).Value = 0
$ws.Cells.Item(206, 8).Value = 0

$ws.Cells.Item(209, 1).Value = "Islas Malvinas"
$ws.Cells.Item(209, 2).Value = 13
$ws.Cells.Item(209, 3).Value = 0
$ws.Cells.Item(209, 4).Value = 13
$ws.Cells.Item(209, 5).Value = 0
$ws.Cells.Item(209, 6).Value = 0
$ws.Cells.Item(209, 7).Value = 0
$ws.Cells.Item(209, 8).Value = 0

$ws.Cells.Item(210, 1).Value = "Groenlandia"
$ws.Cells.Item(210, 2).Value = 13
$ws.Cells.Item(210, 3).Value = 0
$ws.Cells.Item(210, 4).Value = 13
$ws.Cells.Item(210, 5).Value = 0
$ws.Cells.Item(210, 6).Value = 0
$ws.Cells.Item(210, 7).Value = 0
$ws.Cells.Item(210, 8).Value = 0
